$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.462.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.571.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'289.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "'0.3693"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("D8").Value = "'50.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("D9").Value = "'0.3385"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").Value = "'1.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("D11").Value = "'0.07524"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "'21.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "'6.991"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "'1.570.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'0.00001119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").Value = "'90.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'0.06769"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'6.367"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.45%  "
$ws.Range("D22").Value = "'16.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Value = "'12.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("D24").Value = "'22.465.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").Value = "'2.369"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'2.649"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.88%  "
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'149.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "'5.051"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "'124.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").Value = "'1.745.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "'1.060"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.27%  "
$ws.Range("D33").Value = "'6.222"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("D34").Value = "'2.016"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "'9.812"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "'0.08369"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'0.02475"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").Value = "'0.06504"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").Value = "'5.410"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").Value = "'0.6216"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'0.5855"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").Value = "'2.063"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").Value = "'125.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("D51").Value = "'0.07296"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
